$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores runs/balls/fours/sixes as text-typed numbers. We set each
# target cell via a text-literal formula (so the written value is a string,
# matching the existing cell type) and then flatten it to a plain value with
# a values-only paste so no formula or formatting is left behind.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# Row 2: runs 8->2, balls 5->3, fours 0->0, sixes 1->0
Set-TextValue 2 3 "2"
Set-TextValue 2 4 "3"
Set-TextValue 2 6 "0"

# Row 3: runs 31->12, balls 24->6, fours 4->1, sixes 1->1
Set-TextValue 3 3 "12"
Set-TextValue 3 4 "6"
Set-TextValue 3 5 "1"

# Row 4: runs 2->1, balls 3->1, fours 0->0, sixes 0->0
Set-TextValue 4 3 "1"
Set-TextValue 4 4 "1"

# Row 5: runs 10->31, balls 13->24, fours 0->4, sixes 0->1
Set-TextValue 5 3 "31"
Set-TextValue 5 4 "24"
Set-TextValue 5 5 "4"
Set-TextValue 5 6 "1"

# Row 6: runs 1->8, balls 1->5, fours 0->0, sixes 0->1
Set-TextValue 6 3 "8"
Set-TextValue 6 4 "5"
Set-TextValue 6 6 "1"

# Row 7: runs 12->7, balls 6->4, fours 1->1, sixes 1->0
Set-TextValue 7 3 "7"
Set-TextValue 7 4 "4"
Set-TextValue 7 6 "0"

# Row 8: runs 7->10, balls 4->13, fours 1->0, sixes 0->0
Set-TextValue 8 3 "10"
Set-TextValue 8 4 "13"
Set-TextValue 8 5 "0"
